$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column F ("Address") before the existing F ("District") column,
# which shifts the old F column -> G.
$ws.Columns("F").Insert()

# Header
$ws.Range("F2").Value = "Address"

# Per-row address values extracted from each teacher's combined Name/Address cell (column B).
$addresses = @{
    3 = 'G H S HirethogaleriDavanagere South'
    4 = 'G H S Ramagondanahalli'
    5 = 'Sri Kalidasa High School Harihar'
    6 = 'S T J High School'
    7 = 'G H S NitturHarihara'
    8 = 'G H S MudhahadadiDavanagere South'
    9 = 'G H S BelaguttiMalligenahalliHonnali'
    10 = 'Smt. Halamma Shamanur ShivappaHigh School Harihara'
    11 = 'G J C High School SectionSasvehalliHonnali'
    12 = 'S G R K High SchoolHarihar'
    13 = 'Abinaya Bharathi H S Davanagere South'
    14 = 'G H S Ganganakatte'
    15 = 'Sri R G N High School Rampura'
    16 = 'National High School ChiluruHonnali'
    17 = 'Sri Karisiddeswara Resi. High School Shyagale'
    19 = 'Honnali'
    20 = 'S B H SchoolKanagondanahalliDavanagere South'
    21 = 'G H P S KammaragatteHonnali'
    22 = 'Sri Durgambika Comp. High School Hondada Road'
    23 = 'S N R H S NandigudiHarihar'
    24 = 'Head MasterSri Uma Pragathi High School ChinnikatteHonnali'
    25 = 'S J J H S HirekalamataHonnali'
    26 = 'Baba Saheb Ambedkar High SchoolHarihar'
    27 = 'Sri Marulasiddeshwara High SchoolAnagodu'
    28 = 'G H S GopanalDavanagere South'
    29 = 'Govt. P U College(High School Section) Mayakonda'
    30 = 'R H S NalkundaDavanagere South'
    31 = 'Sri Siddeswara High SchoolRangavvanahalli(Kurki)Davanagere South'
    32 = 'Sri Vinayaka High SchoolBenakanahalliHonnali'
    33 = 'Sri Manjunatheshwara High SchoolAvaragere'
    34 = 'Sri Maruthi High School KulagatteHonnali'
    35 = 'Govt. Junior CollegeHonnali'
    36 = 'S S V R H S Gudal'
    37 = 'S M H S DevarabelekereHarihara'
    38 = 'Govt. Girls P U CollegeNyamathi Honnali'
    39 = 'Sri SiddalingeswaraHigh SchoolG BevinahalliHarihara'
    40 = 'Govt. High School ShiramagondanahallyDavanagere South'
    41 = 'Sri Patel Gurubasappa High School BelludiHarihar'
    42 = 'G G H S Gandhimadan Harihar'
    43 = 'South'
    44 = 'Sharana Sangama High School KumbaluruHarihara'
    45 = 'Akkamahadevi Girls High SchoolDavanagere North'
    46 = 'S R H S MalebennurHarihar'
    47 = 'Shree Malagere Veerappa Channappa High SchoolVodeyara HatturuHonnali'
    48 = 'Govt Upgrade Primary SchoolKenchanahalliHarihar'
    49 = 'G H S Huvinamadu'
    50 = 'South'
    51 = 'Govt. High School LingapuraHonnali'
    52 = 'R S G G J C Davanagere South'
    53 = 'G H S NagarasanahallyDavanagere South'
    54 = 'S S High School Hemmana Bethur'
    55 = 'G H S KokkanurHarihar'
    56 = 'Sri Maruthi High SchoolAnjaneyanagara(Honnamara) Kukkuwade'
    57 = 'G H S ShamanurDavanagere South'
    58 = 'M M M G H S K R RoadDavanagere North'
    60 = 'Rural High School Mandalur'
    61 = 'Govt. High School HalivanaHarihar'
    62 = 'A G H S ThurchagattaDavanagere South'
    63 = 'Y B S G High SchoolVinobanagara'
}

foreach ($row in $addresses.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $addresses[$row]
}
